$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-20 (B..F) per the diff ---
$ws.Range("B2").Value = 'NSE:ADSL'
$ws.Range("C2").Value = 'NSE:AARTECH'
$ws.Range("D2").Value = 'NSE:BAJAJ-AUTO'
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 'NSE:CIPLA'

$ws.Range("B3").Value = 'NSE:AJOONI'
$ws.Range("C3").Value = 'NSE:AARTIDRUGS'
$ws.Range("D3").Value = 'NSE:HAVELLS'
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = 'NSE:COLPAL'

$ws.Range("B4").Value = 'NSE:ANDHRSUGAR'
$ws.Range("C4").Value = 'NSE:AGSTRA'
$ws.Range("D4").Value = 'NSE:MCX'
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 'NSE:HINDPETRO'

$ws.Range("B5").Value = 'NSE:BHARTIARTL'
$ws.Range("C5").Value = 'NSE:ALANKIT'
$ws.Range("D5").Value = 'NSE:OBEROIRLTY'
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 'NSE:INDUSTOWER'

$ws.Range("B6").Value = 'NSE:CENTUM'
$ws.Range("C6").Value = 'NSE:ASHAPURMIN'
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()

$ws.Range("B7").Value = 'NSE:CHOICEIN'
$ws.Range("C7").Value = 'NSE:ASIANENE'
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()

$ws.Range("B8").Value = 'NSE:CIGNITITEC'
$ws.Range("C8").Value = 'NSE:ATGL'
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

$ws.Range("B9").Value = 'NSE:COLPAL'
$ws.Range("C9").Value = 'NSE:AWL'
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()

$ws.Range("B10").Value = 'NSE:DBOL'
$ws.Range("C10").Value = 'NSE:BARBEQUE'
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()

$ws.Range("B11").Value = 'NSE:DHAMPURSUG'
$ws.Range("C11").Value = 'NSE:BIGBLOC'
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()

$ws.Range("B12").Value = 'NSE:DODLA'
$ws.Range("C12").Value = 'NSE:BODALCHEM'
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()

$ws.Range("B13").Value = 'NSE:EMUDHRA'
$ws.Range("C13").Value = 'NSE:BSL'
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("F13").ClearContents()

$ws.Range("B14").Value = 'NSE:ETHOSLTD'
$ws.Range("C14").Value = 'NSE:CAMLINFINE'
$ws.Range("D14").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("F14").ClearContents()

$ws.Range("B15").Value = 'NSE:GSPL'
$ws.Range("C15").Value = 'NSE:CASTROLIND'
$ws.Range("D15").ClearContents()
$ws.Range("E15").ClearContents()
$ws.Range("F15").ClearContents()

$ws.Range("B16").Value = 'NSE:HARDWYN'
$ws.Range("C16").Value = 'NSE:CENTURYTEX'
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()

$ws.Range("B17").Value = 'NSE:HIKAL'
$ws.Range("C17").Value = 'NSE:CHEMCON'
$ws.Range("D17").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("F17").ClearContents()

$ws.Range("B18").Value = 'NSE:HINDPETRO'
$ws.Range("C18").Value = 'NSE:DCW'
$ws.Range("D18").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("F18").ClearContents()

$ws.Range("B19").Value = 'NSE:INFRABEES'
$ws.Range("C19").Value = 'NSE:DHANI'
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()

$ws.Range("B20").Value = 'NSE:IOC'
$ws.Range("C20").Value = 'NSE:DIAMONDYD'
$ws.Range("D20").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("F20").ClearContents()

# --- Add new rows 21-43 ---
# Copy the number-column formatting (style) from A20 down to the new A cells first
$ws.Range("A20").Copy()
$ws.Range("A21:A43").PasteSpecial(-4122)

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 'NSE:JLHL'
$ws.Range("C21").Value = 'NSE:DISHTV'
$ws.Range("D21").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("F21").ClearContents()

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 'NSE:NV20BEES'
$ws.Range("C22").Value = 'NSE:DOLATALGO'
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("F22").ClearContents()

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 'NSE:PLASTIBLEN'
$ws.Range("C23").Value = 'NSE:EPIGRAL'
$ws.Range("D23").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("F23").ClearContents()

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 'NSE:PRESTIGE'
$ws.Range("C24").Value = 'NSE:GHCLTEXTIL'
$ws.Range("D24").ClearContents()
$ws.Range("E24").ClearContents()
$ws.Range("F24").ClearContents()

$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 'NSE:PRIMESECU'
$ws.Range("C25").Value = 'NSE:GLAXO'
$ws.Range("D25").ClearContents()
$ws.Range("E25").ClearContents()
$ws.Range("F25").ClearContents()

$ws.Range("A26").Value = 24
$ws.Range("B26").ClearContents()
$ws.Range("C26").Value = 'NSE:IMFA'
$ws.Range("D26").ClearContents()
$ws.Range("E26").ClearContents()
$ws.Range("F26").ClearContents()

$ws.Range("A27").Value = 25
$ws.Range("B27").ClearContents()
$ws.Range("C27").Value = 'NSE:INDNIPPON'
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("F27").ClearContents()

$ws.Range("A28").Value = 26
$ws.Range("B28").ClearContents()
$ws.Range("C28").Value = 'NSE:JPOLYINVST'
$ws.Range("D28").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("F28").ClearContents()

$ws.Range("A29").Value = 27
$ws.Range("B29").ClearContents()
$ws.Range("C29").Value = 'NSE:JPPOWER'
$ws.Range("D29").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("F29").ClearContents()

$ws.Range("A30").Value = 28
$ws.Range("B30").ClearContents()
$ws.Range("C30").Value = 'NSE:JSWINFRA'
$ws.Range("D30").ClearContents()
$ws.Range("E30").ClearContents()
$ws.Range("F30").ClearContents()

$ws.Range("A31").Value = 29
$ws.Range("B31").ClearContents()
$ws.Range("C31").Value = 'NSE:KHAICHEM'
$ws.Range("D31").ClearContents()
$ws.Range("E31").ClearContents()
$ws.Range("F31").ClearContents()

$ws.Range("A32").Value = 30
$ws.Range("B32").ClearContents()
$ws.Range("C32").Value = 'NSE:LUMAXTECH'
$ws.Range("D32").ClearContents()
$ws.Range("E32").ClearContents()
$ws.Range("F32").ClearContents()

$ws.Range("A33").Value = 31
$ws.Range("B33").ClearContents()
$ws.Range("C33").Value = 'NSE:MARALOVER'
$ws.Range("D33").ClearContents()
$ws.Range("E33").ClearContents()
$ws.Range("F33").ClearContents()

$ws.Range("A34").Value = 32
$ws.Range("B34").ClearContents()
$ws.Range("C34").Value = 'NSE:MOL'
$ws.Range("D34").ClearContents()
$ws.Range("E34").ClearContents()
$ws.Range("F34").ClearContents()

$ws.Range("A35").Value = 33
$ws.Range("B35").ClearContents()
$ws.Range("C35").Value = 'NSE:MUNJALAU'
$ws.Range("D35").ClearContents()
$ws.Range("E35").ClearContents()
$ws.Range("F35").ClearContents()

$ws.Range("A36").Value = 34
$ws.Range("B36").ClearContents()
$ws.Range("C36").Value = 'NSE:NYKAA'
$ws.Range("D36").ClearContents()
$ws.Range("E36").ClearContents()
$ws.Range("F36").ClearContents()

$ws.Range("A37").Value = 35
$ws.Range("B37").ClearContents()
$ws.Range("C37").Value = 'NSE:PARAGMILK'
$ws.Range("D37").ClearContents()
$ws.Range("E37").ClearContents()
$ws.Range("F37").ClearContents()

$ws.Range("A38").Value = 36
$ws.Range("B38").ClearContents()
$ws.Range("C38").Value = 'NSE:PRECAM'
$ws.Range("D38").ClearContents()
$ws.Range("E38").ClearContents()
$ws.Range("F38").ClearContents()

$ws.Range("A39").Value = 37
$ws.Range("B39").ClearContents()
$ws.Range("C39").Value = 'NSE:PRIVISCL'
$ws.Range("D39").ClearContents()
$ws.Range("E39").ClearContents()
$ws.Range("F39").ClearContents()

$ws.Range("A40").Value = 38
$ws.Range("B40").ClearContents()
$ws.Range("C40").Value = 'NSE:PTL'
$ws.Range("D40").ClearContents()
$ws.Range("E40").ClearContents()
$ws.Range("F40").ClearContents()

$ws.Range("A41").Value = 39
$ws.Range("B41").ClearContents()
$ws.Range("C41").Value = 'NSE:QUESS'
$ws.Range("D41").ClearContents()
$ws.Range("E41").ClearContents()
$ws.Range("F41").ClearContents()

$ws.Range("A42").Value = 40
$ws.Range("B42").ClearContents()
$ws.Range("C42").Value = 'NSE:RADHIKAJWE'
$ws.Range("D42").ClearContents()
$ws.Range("E42").ClearContents()
$ws.Range("F42").ClearContents()

$ws.Range("A43").Value = 41
$ws.Range("B43").ClearContents()
$ws.Range("C43").Value = 'NSE:RML'
$ws.Range("D43").ClearContents()
$ws.Range("E43").ClearContents()
$ws.Range("F43").ClearContents()

